$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.721.33'
$ws.Range("E2").Value = '  -0.10%  '
$ws.Range("D3").Value = '2.291.92'
$ws.Range("E3").Value = '  -1.11%  '
$ws.Range("E4").Value = '  -0.08%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '96.40'
$c.Style = "Normal"
$ws.Range("E5").Value = '  +4.28%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '268.01'
$c.Style = "Normal"
$ws.Range("E6").Value = '  -0.24%  '
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.619'
$c.Style = "Normal"
$ws.Range("E7").Value = '  -1.54%  '
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.999'
$c.Style = "Normal"
$ws.Range("E8").Value = '  -0.13%  '
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.611'
$c.Style = "Normal"
$ws.Range("E9").Value = '  -1.50%  '
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '45.95'
$c.Style = "Normal"
$ws.Range("E10").Value = '  +2.51%  '
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.0937'
$c.Style = "Normal"
$ws.Range("E11").Value = '  +0.19%  '
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '7.89'
$c.Style = "Normal"
$ws.Range("E12").Value = '  -1.23%  '
$ws.Range("E13").Value = '  +0.28%  '
$ws.Range("D14").Value = '2.632.19'
$ws.Range("E14").Value = '  -1.23%  '
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '15.18'
$c.Style = "Normal"
$ws.Range("E15").Value = '  -0.84%  '
$ws.Range("E16").Value = '  -0.82%  '
$ws.Range("D17").Value = '2.291.94'
$ws.Range("E17").Value = '  -1.45%  '
$ws.Range("D18").Value = '43.591.89'
$ws.Range("E18").Value = '  -0.39%  '
$ws.Range("E19").Value = '  +1.75%  '
$ws.Range("E20").Value = '  -1.52%  '
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '72.28'
$c.Style = "Normal"
$ws.Range("E21").Value = '  +1.40%  '
$ws.Range("E22").Value = '  +9.86%  '
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '233.13'
$c.Style = "Normal"
$ws.Range("E23").Value = '  -3.56%  '
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '9.16'
$c.Style = "Normal"
$ws.Range("E24").Value = '  -5.28%  '
$ws.Range("E25").Value = '  -0.05%  '
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '2.53'
$c.Style = "Normal"
$ws.Range("E26").Value = '  +1.43%  '
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '11.23'
$c.Style = "Normal"
$ws.Range("E27").Value = '  -0.20%  '
$ws.Range("E28").Value = '  +3.03%  '
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '40.22'
$c.Style = "Normal"
$ws.Range("E29").Value = '  +3.24%  '
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '2.23'
$c.Style = "Normal"
$ws.Range("E30").Value = '  -3.81%  '
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '175.60'
$c.Style = "Normal"
$ws.Range("E31").Value = '  +1.79%  '
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '21.87'
$c.Style = "Normal"
$ws.Range("E32").Value = '  -3.13%  '
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '0.0895'
$c.Style = "Normal"
$ws.Range("E33").Value = '  +0.18%  '
$ws.Range("E34").Value = '  -3.21%  '
$ws.Range("E35").Value = '  -0.10%  '
$ws.Range("E36").Value = '  -2.28%  '
$ws.Range("E37").Value = '  +1.82%  '
$ws.Range("E38").Value = '  -3.42%  '
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '3.41'
$c.Style = "Normal"
$ws.Range("E39").Value = '  +1.46%  '
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '0.246'
$c.Style = "Normal"
$ws.Range("E40").Value = '  +3.79%  '
$ws.Range("E41").Value = '  +0.69%  '
$ws.Range("B42").Value = 'ARBITRUM'
$ws.Range("C42").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '1.36'
$c.Style = "Normal"
$ws.Range("E42").Value = '  +2.51%  '
$ws.Range("B43").Value = 'Celestia'
$ws.Range("C43").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '12.28'
$c.Style = "Normal"
$ws.Range("E43").Value = '  +0.78%  '
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '64.73'
$c.Style = "Normal"
$ws.Range("E44").Value = '  +5.62%  '
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '8.81'
$c.Style = "Normal"
$ws.Range("E45").Value = '  -1.23%  '
$ws.Range("E47").Value = '  -0.08%  '
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '97.48'
$c.Style = "Normal"
$ws.Range("E48").Value = '  -2.95%  '
$ws.Range("E49").Value = '  -0.52%  '
$ws.Range("B50").Value = 'TheGraph'
$ws.Range("C50").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '0.186'
$c.Style = "Normal"
$ws.Range("E50").Value = '  +8.79%  '
$ws.Range("B51").Value = 'RocketPoolETH'
$ws.Range("C51").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D51").Value = '2.511.07'
$ws.Range("E51").Value = '  -1.25%  '
